# Auto-generated edit script applying numeric updates described in the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7823.8184
$ws.Range("J40").Value = 8063.231
$ws.Range("L40").Value = 8063.231
$ws.Range("N40").Value = -8413.231
$ws.Range("H62").Value = 8338343.5
$ws.Range("J62").Value = 7059.75
$ws.Range("L62").Value = 7059.75
$ws.Range("N62").Value = -8307.75
$ws.Range("H65").Value = 8338343.5
$ws.Range("J65").Value = 7059.75
$ws.Range("L65").Value = 35298.75
$ws.Range("N65").Value = -41538.75
$ws.Range("H88").Value = 3620.6
$ws.Range("I88").Value = 3369
$ws.Range("J88").Value = 3683.5
$ws.Range("K88").Value = 3369
$ws.Range("L88").Value = 3683.5
$ws.Range("M88").Value = -2963
$ws.Range("N88").Value = -4495.5
$ws.Range("H91").Value = 3620.6
$ws.Range("I91").Value = 3369
$ws.Range("J91").Value = 3683.5
$ws.Range("K91").Value = 3369
$ws.Range("L91").Value = 3683.5
$ws.Range("M91").Value = -1965
$ws.Range("N91").Value = -6491.5
$ws.Range("H116").Value = 17135.555
$ws.Range("J116").Value = 30139.5
$ws.Range("L116").Value = 30139.5
$ws.Range("N116").Value = -37023.5
$ws.Range("H132").Value = 4255
$ws.Range("I132").Value = 4333.615
$ws.Range("K132").Value = 13000.845
$ws.Range("M132").Value = -10470.845
$ws.Range("H137").Value = 627337.25
$ws.Range("I137").Value = 557213.75
$ws.Range("J137").Value = 717495.9399999999
$ws.Range("K137").Value = 1671641.25
$ws.Range("L137").Value = 2152487.82
$ws.Range("M137").Value = -1669091.25
$ws.Range("N137").Value = -2157587.82
$ws.Range("H138").Value = 5454.206
$ws.Range("I138").Value = 1554.6666
$ws.Range("J138").Value = 11753.462
$ws.Range("K138").Value = 4663.9998
$ws.Range("L138").Value = 35260.386
$ws.Range("M138").Value = 476.0002000000004
$ws.Range("N138").Value = -45540.386

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5235.857
$ws.Range("I32").Value = 5235.857
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5235.857
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -4948.857
$ws.Range("N32").Value = ""
$ws.Range("H63").Value = 6941.6665
$ws.Range("J63").Value = 9300
$ws.Range("L63").Value = 9300
$ws.Range("N63").Value = -10672
$ws.Range("H66").Value = 6941.6665
$ws.Range("J66").Value = 9300
$ws.Range("L66").Value = 46500
$ws.Range("N66").Value = -53364

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2470.4211
$ws.Range("I20").Value = 1779
$ws.Range("J20").Value = 3421.125
$ws.Range("K20").Value = 1779
$ws.Range("L20").Value = 3421.125
$ws.Range("M20").Value = -1532
$ws.Range("N20").Value = -3915.125
$ws.Range("H22").Value = 254
$ws.Range("I22").Value = 216.6
$ws.Range("J22").Value = 291.4
$ws.Range("K22").Value = 216.6
$ws.Range("L22").Value = 291.4
$ws.Range("M22").Value = -43.59999999999999
$ws.Range("N22").Value = -637.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27525.1
$ws.Range("I31").Value = 1584.3103
$ws.Range("J31").Value = 95914.45
$ws.Range("K31").Value = 1584.3103
$ws.Range("L31").Value = 95914.45
$ws.Range("M31").Value = -1289.3103
$ws.Range("N31").Value = -96504.45
$ws.Range("H34").Value = 27525.1
$ws.Range("I34").Value = 1584.3103
$ws.Range("J34").Value = 95914.45
$ws.Range("K34").Value = 1584.3103
$ws.Range("L34").Value = 95914.45
$ws.Range("M34").Value = -1382.3103
$ws.Range("N34").Value = -96318.45
$ws.Range("H107").Value = 339.125
$ws.Range("I107").Value = 339.125
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 339.125
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1580.875
$ws.Range("N107").Value = ""
$ws.Range("H134").Value = 802782.4
$ws.Range("I134").Value = 502526.66
$ws.Range("K134").Value = 1507579.98
$ws.Range("M134").Value = -1505044.98

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 500249.5
$ws.Range("I70").Value = 500249.5
$ws.Range("K70").Value = 1500748.5
$ws.Range("M70").Value = -1500433.5
$ws.Range("H73").Value = 500249.5
$ws.Range("I73").Value = 500249.5
$ws.Range("K73").Value = 1500748.5
$ws.Range("M73").Value = -1499656.5
$ws.Range("H122").Value = 42380.707
$ws.Range("J122").Value = 91884.37
$ws.Range("L122").Value = 826959.33
$ws.Range("N122").Value = -831859.33
$ws.Range("H129").Value = 1426.6
$ws.Range("I129").Value = 600
$ws.Range("K129").Value = 1800
$ws.Range("M129").Value = 3200
$ws.Range("H131").Value = 12900927
$ws.Range("I131").Value = 33433984
$ws.Range("J131").Value = 67766.5
$ws.Range("K131").Value = 100301952
$ws.Range("L131").Value = 203299.5
$ws.Range("M131").Value = -100296912
$ws.Range("N131").Value = -213379.5
$ws.Range("H132").Value = 614857.3
$ws.Range("I132").Value = 112642
$ws.Range("J132").Value = 1117072.6
$ws.Range("K132").Value = 1013778
$ws.Range("L132").Value = 10053653.4
$ws.Range("M132").Value = -1011248
$ws.Range("N132").Value = -10058713.4
$ws.Range("H137").Value = 2150.6365
$ws.Range("J137").Value = 10000
$ws.Range("L137").Value = 30000
$ws.Range("N137").Value = -40200

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4246.857
$ws.Range("I122").Value = 2287.7144
$ws.Range("J122").Value = 6206
$ws.Range("K122").Value = 6863.1432
$ws.Range("L122").Value = 18618
$ws.Range("M122").Value = -4413.1432
$ws.Range("N122").Value = -23518

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 6842.2
$ws.Range("J74").Value = 7366.875
$ws.Range("L74").Value = 7366.875
$ws.Range("N74").Value = -9238.875
$ws.Range("H77").Value = 6842.2
$ws.Range("J77").Value = 7366.875
$ws.Range("L77").Value = 22100.625
$ws.Range("N77").Value = -31460.625
$ws.Range("H107").Value = 499.42105
$ws.Range("I107").Value = 564.9167
$ws.Range("K107").Value = 1694.7501
$ws.Range("M107").Value = 225.2499
$ws.Range("H135").Value = 69799.8
$ws.Range("J135").Value = 69799.8
$ws.Range("L135").Value = 69799.8
$ws.Range("N135").Value = -79939.8
$ws.Range("H136").Value = 7675353.5
$ws.Range("I136").Value = 8374901
$ws.Range("J136").Value = 504993.5
$ws.Range("K136").Value = 25124703
$ws.Range("L136").Value = 1514980.5
$ws.Range("M136").Value = -25122153
$ws.Range("N136").Value = -1520080.5
